$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 190-213 (new sensor log entries) ---
$ws = $wb.Worksheets.Item("PIR")
$ws.Cells.Item(190, 1).NumberFormat = "@"
$ws.Cells.Item(190, 1).Value = '2026-02-01'
$ws.Cells.Item(190, 2).Value = '14:06:27'
$ws.Cells.Item(190, 3).Value = '14:00'
$ws.Cells.Item(190, 4).Value = 'Bathroom'
$ws.Cells.Item(190, 5).Value = 'No Motion'
$ws.Cells.Item(190, 6).Value = 'Inactive'
$ws.Cells.Item(191, 1).NumberFormat = "@"
$ws.Cells.Item(191, 1).Value = '2026-02-01'
$ws.Cells.Item(191, 2).Value = '14:06:28'
$ws.Cells.Item(191, 3).Value = '14:00'
$ws.Cells.Item(191, 4).Value = 'Bathroom'
$ws.Cells.Item(191, 5).Value = 'No Motion'
$ws.Cells.Item(191, 6).Value = 'Inactive'
$ws.Cells.Item(192, 1).NumberFormat = "@"
$ws.Cells.Item(192, 1).Value = '2026-02-01'
$ws.Cells.Item(192, 2).Value = '14:06:28'
$ws.Cells.Item(192, 3).Value = '14:00'
$ws.Cells.Item(192, 4).Value = 'Bathroom'
$ws.Cells.Item(192, 5).Value = 'No Motion'
$ws.Cells.Item(192, 6).Value = 'Inactive'
$ws.Cells.Item(193, 1).NumberFormat = "@"
$ws.Cells.Item(193, 1).Value = '2026-02-01'
$ws.Cells.Item(193, 2).Value = '14:06:29'
$ws.Cells.Item(193, 3).Value = '14:00'
$ws.Cells.Item(193, 4).Value = 'Bathroom'
$ws.Cells.Item(193, 5).Value = 'No Motion'
$ws.Cells.Item(193, 6).Value = 'Inactive'
$ws.Cells.Item(194, 1).NumberFormat = "@"
$ws.Cells.Item(194, 1).Value = '2026-02-01'
$ws.Cells.Item(194, 2).Value = '14:06:33'
$ws.Cells.Item(194, 3).Value = '14:00'
$ws.Cells.Item(194, 4).Value = 'Bathroom'
$ws.Cells.Item(194, 5).Value = 'No Motion'
$ws.Cells.Item(194, 6).Value = 'Inactive'
$ws.Cells.Item(195, 1).NumberFormat = "@"
$ws.Cells.Item(195, 1).Value = '2026-02-01'
$ws.Cells.Item(195, 2).Value = '14:06:33'
$ws.Cells.Item(195, 3).Value = '14:00'
$ws.Cells.Item(195, 4).Value = 'Bathroom'
$ws.Cells.Item(195, 5).Value = 'No Motion'
$ws.Cells.Item(195, 6).Value = 'Inactive'
$ws.Cells.Item(196, 1).NumberFormat = "@"
$ws.Cells.Item(196, 1).Value = '2026-02-01'
$ws.Cells.Item(196, 2).Value = '14:06:38'
$ws.Cells.Item(196, 3).Value = '14:00'
$ws.Cells.Item(196, 4).Value = 'Bathroom'
$ws.Cells.Item(196, 5).Value = 'No Motion'
$ws.Cells.Item(196, 6).Value = 'Inactive'
$ws.Cells.Item(197, 1).NumberFormat = "@"
$ws.Cells.Item(197, 1).Value = '2026-02-01'
$ws.Cells.Item(197, 2).Value = '14:06:38'
$ws.Cells.Item(197, 3).Value = '14:00'
$ws.Cells.Item(197, 4).Value = 'Bathroom'
$ws.Cells.Item(197, 5).Value = 'No Motion'
$ws.Cells.Item(197, 6).Value = 'Inactive'
$ws.Cells.Item(198, 1).NumberFormat = "@"
$ws.Cells.Item(198, 1).Value = '2026-02-01'
$ws.Cells.Item(198, 2).Value = '14:06:43'
$ws.Cells.Item(198, 3).Value = '14:00'
$ws.Cells.Item(198, 4).Value = 'Bathroom'
$ws.Cells.Item(198, 5).Value = 'No Motion'
$ws.Cells.Item(198, 6).Value = 'Inactive'
$ws.Cells.Item(199, 1).NumberFormat = "@"
$ws.Cells.Item(199, 1).Value = '2026-02-01'
$ws.Cells.Item(199, 2).Value = '14:06:43'
$ws.Cells.Item(199, 3).Value = '14:00'
$ws.Cells.Item(199, 4).Value = 'Bathroom'
$ws.Cells.Item(199, 5).Value = 'No Motion'
$ws.Cells.Item(199, 6).Value = 'Inactive'
$ws.Cells.Item(200, 1).NumberFormat = "@"
$ws.Cells.Item(200, 1).Value = '2026-02-01'
$ws.Cells.Item(200, 2).Value = '14:06:48'
$ws.Cells.Item(200, 3).Value = '14:00'
$ws.Cells.Item(200, 4).Value = 'Bathroom'
$ws.Cells.Item(200, 5).Value = 'No Motion'
$ws.Cells.Item(200, 6).Value = 'Inactive'
$ws.Cells.Item(201, 1).NumberFormat = "@"
$ws.Cells.Item(201, 1).Value = '2026-02-01'
$ws.Cells.Item(201, 2).Value = '14:06:48'
$ws.Cells.Item(201, 3).Value = '14:00'
$ws.Cells.Item(201, 4).Value = 'Bathroom'
$ws.Cells.Item(201, 5).Value = 'No Motion'
$ws.Cells.Item(201, 6).Value = 'Inactive'
$ws.Cells.Item(202, 1).NumberFormat = "@"
$ws.Cells.Item(202, 1).Value = '2026-02-01'
$ws.Cells.Item(202, 2).Value = '14:10:29'
$ws.Cells.Item(202, 3).Value = '14:00'
$ws.Cells.Item(202, 4).Value = 'Bathroom'
$ws.Cells.Item(202, 5).Value = 'No Motion'
$ws.Cells.Item(202, 6).Value = 'Inactive'
$ws.Cells.Item(203, 1).NumberFormat = "@"
$ws.Cells.Item(203, 1).Value = '2026-02-01'
$ws.Cells.Item(203, 2).Value = '14:10:34'
$ws.Cells.Item(203, 3).Value = '14:00'
$ws.Cells.Item(203, 4).Value = 'Bathroom'
$ws.Cells.Item(203, 5).Value = 'No Motion'
$ws.Cells.Item(203, 6).Value = 'Inactive'
$ws.Cells.Item(204, 1).NumberFormat = "@"
$ws.Cells.Item(204, 1).Value = '2026-02-01'
$ws.Cells.Item(204, 2).Value = '14:10:39'
$ws.Cells.Item(204, 3).Value = '14:00'
$ws.Cells.Item(204, 4).Value = 'Bathroom'
$ws.Cells.Item(204, 5).Value = 'No Motion'
$ws.Cells.Item(204, 6).Value = 'Inactive'
$ws.Cells.Item(205, 1).NumberFormat = "@"
$ws.Cells.Item(205, 1).Value = '2026-02-01'
$ws.Cells.Item(205, 2).Value = '14:10:50'
$ws.Cells.Item(205, 3).Value = '14:00'
$ws.Cells.Item(205, 4).Value = 'Bathroom'
$ws.Cells.Item(205, 5).Value = 'Motion Detected'
$ws.Cells.Item(205, 6).Value = 'Active'
$ws.Cells.Item(206, 1).NumberFormat = "@"
$ws.Cells.Item(206, 1).Value = '2026-02-01'
$ws.Cells.Item(206, 2).Value = '14:10:52'
$ws.Cells.Item(206, 3).Value = '14:00'
$ws.Cells.Item(206, 4).Value = 'Bathroom'
$ws.Cells.Item(206, 5).Value = 'No Motion'
$ws.Cells.Item(206, 6).Value = 'Inactive'
$ws.Cells.Item(207, 1).NumberFormat = "@"
$ws.Cells.Item(207, 1).Value = '2026-02-01'
$ws.Cells.Item(207, 2).Value = '14:10:52'
$ws.Cells.Item(207, 3).Value = '14:00'
$ws.Cells.Item(207, 4).Value = 'Bathroom'
$ws.Cells.Item(207, 5).Value = 'Motion Detected'
$ws.Cells.Item(207, 6).Value = 'Active'
$ws.Cells.Item(208, 1).NumberFormat = "@"
$ws.Cells.Item(208, 1).Value = '2026-02-01'
$ws.Cells.Item(208, 2).Value = '14:10:57'
$ws.Cells.Item(208, 3).Value = '14:00'
$ws.Cells.Item(208, 4).Value = 'Bathroom'
$ws.Cells.Item(208, 5).Value = 'No Motion'
$ws.Cells.Item(208, 6).Value = 'Inactive'
$ws.Cells.Item(209, 1).NumberFormat = "@"
$ws.Cells.Item(209, 1).Value = '2026-02-01'
$ws.Cells.Item(209, 2).Value = '14:11:02'
$ws.Cells.Item(209, 3).Value = '14:00'
$ws.Cells.Item(209, 4).Value = 'Bathroom'
$ws.Cells.Item(209, 5).Value = 'No Motion'
$ws.Cells.Item(209, 6).Value = 'Inactive'
$ws.Cells.Item(210, 1).NumberFormat = "@"
$ws.Cells.Item(210, 1).Value = '2026-02-01'
$ws.Cells.Item(210, 2).Value = '14:11:07'
$ws.Cells.Item(210, 3).Value = '14:00'
$ws.Cells.Item(210, 4).Value = 'Bathroom'
$ws.Cells.Item(210, 5).Value = 'No Motion'
$ws.Cells.Item(210, 6).Value = 'Inactive'
$ws.Cells.Item(211, 1).NumberFormat = "@"
$ws.Cells.Item(211, 1).Value = '2026-02-01'
$ws.Cells.Item(211, 2).Value = '14:11:12'
$ws.Cells.Item(211, 3).Value = '14:00'
$ws.Cells.Item(211, 4).Value = 'Bathroom'
$ws.Cells.Item(211, 5).Value = 'No Motion'
$ws.Cells.Item(211, 6).Value = 'Inactive'
$ws.Cells.Item(212, 1).NumberFormat = "@"
$ws.Cells.Item(212, 1).Value = '2026-02-01'
$ws.Cells.Item(212, 2).Value = '14:11:15'
$ws.Cells.Item(212, 3).Value = '14:00'
$ws.Cells.Item(212, 4).Value = 'Bathroom'
$ws.Cells.Item(212, 5).Value = 'Motion Detected'
$ws.Cells.Item(212, 6).Value = 'Active'
$ws.Cells.Item(213, 1).NumberFormat = "@"
$ws.Cells.Item(213, 1).Value = '2026-02-01'
$ws.Cells.Item(213, 2).Value = '14:11:23'
$ws.Cells.Item(213, 3).Value = '14:00'
$ws.Cells.Item(213, 4).Value = 'Bathroom'
$ws.Cells.Item(213, 5).Value = 'No Motion'
$ws.Cells.Item(213, 6).Value = 'Inactive'

# --- Humidity sheet: append rows 120-133 (new sensor log entries) ---
$ws = $wb.Worksheets.Item("Humidity")
$ws.Cells.Item(120, 1).NumberFormat = "@"
$ws.Cells.Item(120, 1).Value = '2026-02-01'
$ws.Cells.Item(120, 2).Value = '14:06:31'
$ws.Cells.Item(120, 3).Value = '14:00'
$ws.Cells.Item(120, 4).Value = 'Bathroom'
$ws.Cells.Item(120, 5).NumberFormat = "@"
$ws.Cells.Item(120, 5).Value = '79.1%'
$ws.Cells.Item(120, 6).Value = 'Active'
$ws.Cells.Item(121, 1).NumberFormat = "@"
$ws.Cells.Item(121, 1).Value = '2026-02-01'
$ws.Cells.Item(121, 2).Value = '14:06:36'
$ws.Cells.Item(121, 3).Value = '14:00'
$ws.Cells.Item(121, 4).Value = 'Bathroom'
$ws.Cells.Item(121, 5).NumberFormat = "@"
$ws.Cells.Item(121, 5).Value = '78.8%'
$ws.Cells.Item(121, 6).Value = 'Active'
$ws.Cells.Item(122, 1).NumberFormat = "@"
$ws.Cells.Item(122, 1).Value = '2026-02-01'
$ws.Cells.Item(122, 2).Value = '14:06:41'
$ws.Cells.Item(122, 3).Value = '14:00'
$ws.Cells.Item(122, 4).Value = 'Bathroom'
$ws.Cells.Item(122, 5).NumberFormat = "@"
$ws.Cells.Item(122, 5).Value = '78.7%'
$ws.Cells.Item(122, 6).Value = 'Active'
$ws.Cells.Item(123, 1).NumberFormat = "@"
$ws.Cells.Item(123, 1).Value = '2026-02-01'
$ws.Cells.Item(123, 2).Value = '14:06:44'
$ws.Cells.Item(123, 3).Value = '14:00'
$ws.Cells.Item(123, 4).Value = 'Bathroom'
$ws.Cells.Item(123, 5).NumberFormat = "@"
$ws.Cells.Item(123, 5).Value = '78.5%'
$ws.Cells.Item(123, 6).Value = 'Active'
$ws.Cells.Item(124, 1).NumberFormat = "@"
$ws.Cells.Item(124, 1).Value = '2026-02-01'
$ws.Cells.Item(124, 2).Value = '14:06:46'
$ws.Cells.Item(124, 3).Value = '14:00'
$ws.Cells.Item(124, 4).Value = 'Bathroom'
$ws.Cells.Item(124, 5).NumberFormat = "@"
$ws.Cells.Item(124, 5).Value = '77.5%'
$ws.Cells.Item(124, 6).Value = 'Active'
$ws.Cells.Item(125, 1).NumberFormat = "@"
$ws.Cells.Item(125, 1).Value = '2026-02-01'
$ws.Cells.Item(125, 2).Value = '14:06:49'
$ws.Cells.Item(125, 3).Value = '14:00'
$ws.Cells.Item(125, 4).Value = 'Bathroom'
$ws.Cells.Item(125, 5).NumberFormat = "@"
$ws.Cells.Item(125, 5).Value = '78.4%'
$ws.Cells.Item(125, 6).Value = 'Active'
$ws.Cells.Item(126, 1).NumberFormat = "@"
$ws.Cells.Item(126, 1).Value = '2026-02-01'
$ws.Cells.Item(126, 2).Value = '14:10:30'
$ws.Cells.Item(126, 3).Value = '14:00'
$ws.Cells.Item(126, 4).Value = 'Bathroom'
$ws.Cells.Item(126, 5).NumberFormat = "@"
$ws.Cells.Item(126, 5).Value = '77.1%'
$ws.Cells.Item(126, 6).Value = 'Active'
$ws.Cells.Item(127, 1).NumberFormat = "@"
$ws.Cells.Item(127, 1).Value = '2026-02-01'
$ws.Cells.Item(127, 2).Value = '14:10:35'
$ws.Cells.Item(127, 3).Value = '14:00'
$ws.Cells.Item(127, 4).Value = 'Bathroom'
$ws.Cells.Item(127, 5).NumberFormat = "@"
$ws.Cells.Item(127, 5).Value = '78.0%'
$ws.Cells.Item(127, 6).Value = 'Active'
$ws.Cells.Item(128, 1).NumberFormat = "@"
$ws.Cells.Item(128, 1).Value = '2026-02-01'
$ws.Cells.Item(128, 2).Value = '14:10:51'
$ws.Cells.Item(128, 3).Value = '14:00'
$ws.Cells.Item(128, 4).Value = 'Bathroom'
$ws.Cells.Item(128, 5).NumberFormat = "@"
$ws.Cells.Item(128, 5).Value = '77.1%'
$ws.Cells.Item(128, 6).Value = 'Active'
$ws.Cells.Item(129, 1).NumberFormat = "@"
$ws.Cells.Item(129, 1).Value = '2026-02-01'
$ws.Cells.Item(129, 2).Value = '14:10:55'
$ws.Cells.Item(129, 3).Value = '14:00'
$ws.Cells.Item(129, 4).Value = 'Bathroom'
$ws.Cells.Item(129, 5).NumberFormat = "@"
$ws.Cells.Item(129, 5).Value = '77.1%'
$ws.Cells.Item(129, 6).Value = 'Active'
$ws.Cells.Item(130, 1).NumberFormat = "@"
$ws.Cells.Item(130, 1).Value = '2026-02-01'
$ws.Cells.Item(130, 2).Value = '14:11:05'
$ws.Cells.Item(130, 3).Value = '14:00'
$ws.Cells.Item(130, 4).Value = 'Bathroom'
$ws.Cells.Item(130, 5).NumberFormat = "@"
$ws.Cells.Item(130, 5).Value = '78.1%'
$ws.Cells.Item(130, 6).Value = 'Active'
$ws.Cells.Item(131, 1).NumberFormat = "@"
$ws.Cells.Item(131, 1).Value = '2026-02-01'
$ws.Cells.Item(131, 2).Value = '14:11:10'
$ws.Cells.Item(131, 3).Value = '14:00'
$ws.Cells.Item(131, 4).Value = 'Bathroom'
$ws.Cells.Item(131, 5).NumberFormat = "@"
$ws.Cells.Item(131, 5).Value = '76.7%'
$ws.Cells.Item(131, 6).Value = 'Active'
$ws.Cells.Item(132, 1).NumberFormat = "@"
$ws.Cells.Item(132, 1).Value = '2026-02-01'
$ws.Cells.Item(132, 2).Value = '14:11:20'
$ws.Cells.Item(132, 3).Value = '14:00'
$ws.Cells.Item(132, 4).Value = 'Bathroom'
$ws.Cells.Item(132, 5).NumberFormat = "@"
$ws.Cells.Item(132, 5).Value = '78.2%'
$ws.Cells.Item(132, 6).Value = 'Active'
$ws.Cells.Item(133, 1).NumberFormat = "@"
$ws.Cells.Item(133, 1).Value = '2026-02-01'
$ws.Cells.Item(133, 2).Value = '14:11:25'
$ws.Cells.Item(133, 3).Value = '14:00'
$ws.Cells.Item(133, 4).Value = 'Bathroom'
$ws.Cells.Item(133, 5).NumberFormat = "@"
$ws.Cells.Item(133, 5).Value = '78.2%'
$ws.Cells.Item(133, 6).Value = 'Active'

# --- Temperature sheet: append rows 41-54 (new sensor log entries) ---
$ws = $wb.Worksheets.Item("Temperature")
$ws.Cells.Item(41, 1).NumberFormat = "@"
$ws.Cells.Item(41, 1).Value = '2026-02-01'
$ws.Cells.Item(41, 2).Value = '14:06:32'
$ws.Cells.Item(41, 3).Value = '14:00'
$ws.Cells.Item(41, 4).Value = 'Bathroom'
$ws.Cells.Item(41, 5).Value = '29.6C'
$ws.Cells.Item(41, 6).Value = 'Active'
$ws.Cells.Item(42, 1).NumberFormat = "@"
$ws.Cells.Item(42, 1).Value = '2026-02-01'
$ws.Cells.Item(42, 2).Value = '14:06:37'
$ws.Cells.Item(42, 3).Value = '14:00'
$ws.Cells.Item(42, 4).Value = 'Bathroom'
$ws.Cells.Item(42, 5).Value = '29.6C'
$ws.Cells.Item(42, 6).Value = 'Active'
$ws.Cells.Item(43, 1).NumberFormat = "@"
$ws.Cells.Item(43, 1).Value = '2026-02-01'
$ws.Cells.Item(43, 2).Value = '14:06:42'
$ws.Cells.Item(43, 3).Value = '14:00'
$ws.Cells.Item(43, 4).Value = 'Bathroom'
$ws.Cells.Item(43, 5).Value = '29.6C'
$ws.Cells.Item(43, 6).Value = 'Active'
$ws.Cells.Item(44, 1).NumberFormat = "@"
$ws.Cells.Item(44, 1).Value = '2026-02-01'
$ws.Cells.Item(44, 2).Value = '14:06:44'
$ws.Cells.Item(44, 3).Value = '14:00'
$ws.Cells.Item(44, 4).Value = 'Bathroom'
$ws.Cells.Item(44, 5).Value = '29.6C'
$ws.Cells.Item(44, 6).Value = 'Active'
$ws.Cells.Item(45, 1).NumberFormat = "@"
$ws.Cells.Item(45, 1).Value = '2026-02-01'
$ws.Cells.Item(45, 2).Value = '14:06:47'
$ws.Cells.Item(45, 3).Value = '14:00'
$ws.Cells.Item(45, 4).Value = 'Bathroom'
$ws.Cells.Item(45, 5).Value = '29.6C'
$ws.Cells.Item(45, 6).Value = 'Active'
$ws.Cells.Item(46, 1).NumberFormat = "@"
$ws.Cells.Item(46, 1).Value = '2026-02-01'
$ws.Cells.Item(46, 2).Value = '14:06:49'
$ws.Cells.Item(46, 3).Value = '14:00'
$ws.Cells.Item(46, 4).Value = 'Bathroom'
$ws.Cells.Item(46, 5).Value = '29.6C'
$ws.Cells.Item(46, 6).Value = 'Active'
$ws.Cells.Item(47, 1).NumberFormat = "@"
$ws.Cells.Item(47, 1).Value = '2026-02-01'
$ws.Cells.Item(47, 2).Value = '14:10:31'
$ws.Cells.Item(47, 3).Value = '14:00'
$ws.Cells.Item(47, 4).Value = 'Bathroom'
$ws.Cells.Item(47, 5).Value = '29.4C'
$ws.Cells.Item(47, 6).Value = 'Active'
$ws.Cells.Item(48, 1).NumberFormat = "@"
$ws.Cells.Item(48, 1).Value = '2026-02-01'
$ws.Cells.Item(48, 2).Value = '14:10:36'
$ws.Cells.Item(48, 3).Value = '14:00'
$ws.Cells.Item(48, 4).Value = 'Bathroom'
$ws.Cells.Item(48, 5).Value = '29.4C'
$ws.Cells.Item(48, 6).Value = 'Active'
$ws.Cells.Item(49, 1).NumberFormat = "@"
$ws.Cells.Item(49, 1).Value = '2026-02-01'
$ws.Cells.Item(49, 2).Value = '14:10:51'
$ws.Cells.Item(49, 3).Value = '14:00'
$ws.Cells.Item(49, 4).Value = 'Bathroom'
$ws.Cells.Item(49, 5).Value = '29.4C'
$ws.Cells.Item(49, 6).Value = 'Active'
$ws.Cells.Item(50, 1).NumberFormat = "@"
$ws.Cells.Item(50, 1).Value = '2026-02-01'
$ws.Cells.Item(50, 2).Value = '14:10:56'
$ws.Cells.Item(50, 3).Value = '14:00'
$ws.Cells.Item(50, 4).Value = 'Bathroom'
$ws.Cells.Item(50, 5).Value = '29.4C'
$ws.Cells.Item(50, 6).Value = 'Active'
$ws.Cells.Item(51, 1).NumberFormat = "@"
$ws.Cells.Item(51, 1).Value = '2026-02-01'
$ws.Cells.Item(51, 2).Value = '14:11:06'
$ws.Cells.Item(51, 3).Value = '14:00'
$ws.Cells.Item(51, 4).Value = 'Bathroom'
$ws.Cells.Item(51, 5).Value = '29.4C'
$ws.Cells.Item(51, 6).Value = 'Active'
$ws.Cells.Item(52, 1).NumberFormat = "@"
$ws.Cells.Item(52, 1).Value = '2026-02-01'
$ws.Cells.Item(52, 2).Value = '14:11:11'
$ws.Cells.Item(52, 3).Value = '14:00'
$ws.Cells.Item(52, 4).Value = 'Bathroom'
$ws.Cells.Item(52, 5).Value = '29.4C'
$ws.Cells.Item(52, 6).Value = 'Active'
$ws.Cells.Item(53, 1).NumberFormat = "@"
$ws.Cells.Item(53, 1).Value = '2026-02-01'
$ws.Cells.Item(53, 2).Value = '14:11:21'
$ws.Cells.Item(53, 3).Value = '14:00'
$ws.Cells.Item(53, 4).Value = 'Bathroom'
$ws.Cells.Item(53, 5).Value = '29.4C'
$ws.Cells.Item(53, 6).Value = 'Active'
$ws.Cells.Item(54, 1).NumberFormat = "@"
$ws.Cells.Item(54, 1).Value = '2026-02-01'
$ws.Cells.Item(54, 2).Value = '14:11:26'
$ws.Cells.Item(54, 3).Value = '14:00'
$ws.Cells.Item(54, 4).Value = 'Bathroom'
$ws.Cells.Item(54, 5).Value = '29.3C'
$ws.Cells.Item(54, 6).Value = 'Active'

# --- Proximity sheet: append rows 25-26 (new sensor log entries) ---
$ws = $wb.Worksheets.Item("Proximity")
$ws.Cells.Item(25, 1).NumberFormat = "@"
$ws.Cells.Item(25, 1).Value = '2026-02-01'
$ws.Cells.Item(25, 2).Value = '14:10:39'
$ws.Cells.Item(25, 3).Value = '14:00'
$ws.Cells.Item(25, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(25, 5).Value = 'ENTER'
$ws.Cells.Item(25, 6).Value = 'User ENTERED Living Room Main Door'
$ws.Cells.Item(26, 1).NumberFormat = "@"
$ws.Cells.Item(26, 1).Value = '2026-02-01'
$ws.Cells.Item(26, 2).Value = '14:10:50'
$ws.Cells.Item(26, 3).Value = '14:00'
$ws.Cells.Item(26, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(26, 5).Value = 'EXIT'
$ws.Cells.Item(26, 6).Value = 'User EXITED Living Room Main Door'

# --- Camera sheet: append rows 15-15 (new sensor log entries) ---
$ws = $wb.Worksheets.Item("Camera")
$ws.Cells.Item(15, 1).NumberFormat = "@"
$ws.Cells.Item(15, 1).Value = '2026-02-01'
$ws.Cells.Item(15, 2).Value = '14:10:49'
$ws.Cells.Item(15, 3).Value = '14:00'
$ws.Cells.Item(15, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(15, 5).Value = 'Image Captured'
$ws.Cells.Item(15, 6).Value = 'Active'
